$wb = $excel.ActiveWorkbook

# --- Step 1: locate the existing sheets we need as references/templates ---
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")

# --- Step 2: insert a new "2022-Q1" sheet right before "总计", using "2021-Q4"
#     as a formatting template (same column layout: fund code/name/size/position/...) ---
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

$q4Sheet.Range("A1:H3").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B-G on data rows hold text (e.g. fund codes keep leading zeros,
# and figures like "5.50" are stored as literal text), matching the other
# quarter sheets. Force text formatting before assigning so Excel doesn't
# coerce numeric-looking strings into numbers.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "014036"
$newSheet.Range("C2").Value = "博时成长回报混合A"
$newSheet.Range("D2").Value = "5.50"
$newSheet.Range("E2").Value = "68.14"
$newSheet.Range("F2").Value = "2.02"
$newSheet.Range("G2").Value = "0.1111"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "014037"
$newSheet.Range("C3").Value = "博时成长回报混合C"
$newSheet.Range("D3").Value = "1.00"
$newSheet.Range("E3").Value = "68.14"
$newSheet.Range("F3").Value = "2.02"
$newSheet.Range("G3").Value = "0.0202"
$newSheet.Range("H3").Value = 10

# --- Step 3: update the "总计" sheet: insert a new top data row with the
#     2022-Q1 summary and renumber the index column. Re-fetch the sheet by
#     name since its position shifted when the new sheet was inserted. ---
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.13

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
